$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new delivery-zone rows for cities whose names start with a number
# (previously mis-parsed because of leading digits in the city name).

# New city/town names first (matches shared-string insertion order).
$ws.Range("B8").Value = "Горки-46"
$ws.Range("B9").Value = "23 Мая"
$ws.Range("G8").Value = "2,1"
$ws.Range("A10").Value = "23 область"
$ws.Range("B10").Value = "Чемихино"

$ws.Range("A8").Value = "Ленинградская обл."
$ws.Range("C8").Value = "2-3"
$ws.Range("D8").Value = "1-2"
$ws.Range("E8").Value = "1-1"
$ws.Range("F8").Value = "47"

$ws.Range("A9").Value = "Ленинградская обл."
$ws.Range("C9").Value = "2-3"
$ws.Range("D9").Value = "1-2"
$ws.Range("E9").Value = "1-1"
$ws.Range("F9").Value = "47"
$ws.Range("G9").Value = "2,1"

$ws.Range("C10").Value = "2-3"
$ws.Range("D10").Value = "1-2"
$ws.Range("E10").Value = "1-1"
$ws.Range("F10").Value = "47"
$ws.Range("G10").Value = "2,1"

$ws.Range("C10:G10").Select()
